# Update the data values in sheet1 (rows 2-21, columns B-G) to reflect
# the refreshed/regenerated report numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @(167, 162, 1500, 1000, 0, 0)
    3  = @(232, 172, 1500, 1000, 27000, 46000)
    4  = @(118, 230, 123, 267, 28536, 45924)
    5  = @(233, 223, 258, 211, 29670, 46631)
    6  = @(235, 178, 137, 213, 31373, 47499)
    7  = @(112, 144, 144, 269, 33696, 47882)
    8  = @(135, 112, 309, 339, 31827, 48477)
    9  = @(206, 145, 259, 467, 30821, 37827)
    10 = @(183, 103, 164, 260, 33620, 37700)
    11 = @(190, 247, 198, 369, 36036, 38007)
    12 = @(244, 191, 202, 161, 38380, 39767)
    13 = @(124, 103, 173, 224, 41693, 42784)
    14 = @(143, 208, 336, 452, 38640, 41132)
    15 = @(200, 114, 289, 197, 41038, 40976)
    16 = @(113, 106, 225, 381, 44775, 43434)
    17 = @(226, 234, 427, 426, 46543, 42174)
    18 = @(177, 189, 212, 194, 40492, 42292)
    19 = @(238, 198, 251, 226, 39658, 41132)
    20 = @(112, 240, 183, 224, 43005, 43904)
    21 = @(238, 119, 414, 186, 39330, 44640)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2   # Column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
